# Adds two new sheets ("Emprestimos" and "Seguros e Cartoes") between
# "Qualidade Cart 2682" and "Margem Financeira", populates them with
# quarterly data, and updates the sheet selections / active tab to match.

$wb = $excel.ActiveWorkbook

$qc = $wb.Worksheets.Item("Qualidade Cart 2682")

# --- New sheet: Emprestimos -------------------------------------------------
$emprestimos = $wb.Worksheets.Add($null, $qc)
$emprestimos.Name = "Emprestimos"

$emprestimos.Range("D4").Value = "3T24"
$emprestimos.Range("E4").Value = "2T25"
$emprestimos.Range("F4").Value = "3T25"

$emprestimos.Range("C5").Value = "Demais"
$emprestimos.Range("D5").Value = 0.139
$emprestimos.Range("E5").Value = 0.119
$emprestimos.Range("F5").Value = 0.093

$emprestimos.Range("C6").Value = "Placas Solares"
$emprestimos.Range("D6").Value = 4.265
$emprestimos.Range("E6").Value = 3.914
$emprestimos.Range("F6").Value = 3.795

$emprestimos.Range("C7").Value = "Consignado Privado"
$emprestimos.Range("D7").Value = 0.525
$emprestimos.Range("E7").Value = 0.457
$emprestimos.Range("F7").Value = 0.395

$emprestimos.Range("C8").Value = "EGV"
$emprestimos.Range("D8").Value = 3.83
$emprestimos.Range("E8").Value = 4.512
$emprestimos.Range("F8").Value = 4.797

$emprestimos.Range("C9").Value = "Total"
$emprestimos.Range("D9").Value = 8.76
$emprestimos.Range("E9").Value = 9.003
$emprestimos.Range("F9").Value = 9.08

# --- New sheet: Seguros e Cartoes -------------------------------------------
$seguros = $wb.Worksheets.Add($null, $emprestimos)
$seguros.Name = "Seguros e Cartoes"

$seguros.Range("C14").Value = "Cartoes"
$seguros.Range("D14").Value = "3T24"
$seguros.Range("E14").Value = "2T25"
$seguros.Range("F14").Value = "3T25"

$seguros.Range("C15").Value = "Total"
$seguros.Range("D15").Value = 4.542
$seguros.Range("E15").Value = 4.811
$seguros.Range("F15").Value = 4.829

# --- Selections / active tab -------------------------------------------------
# Replay the selection changes in the order that leaves "Seguros e Cartoes"
# as the final active sheet/tab (matching tabSelected + bookViews.activeTab).
$emprestimos.Activate() | Out-Null
$emprestimos.Range("F6").Select() | Out-Null

$qc.Activate() | Out-Null
$qc.Range("I22").Select() | Out-Null

$seguros.Activate() | Out-Null
$seguros.Range("E12").Select() | Out-Null

Write-Host "Added Emprestimos and Seguros e Cartoes sheets"
